$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 3402.533333333333, 3671, 3175, 0.05159952640533447),
    @(1, 3593.833333333333, 3855, 3257, 0.05018406709035238),
    @(2, 3688.033333333333, 3963, 3381, 0.05354119936625163),
    @(3, 3503.7, 3755, 3143, 0.05186223189036052),
    @(4, 2647.833333333333, 2933, 2307, 0.05531125068664551),
    @(5, 2830.033333333333, 3060, 2460, 0.05243798891703288),
    @(6, 3603.9, 3888, 3214, 0.05511205196380616),
    @(7, 3156.3, 3438, 2789, 0.05429483254750569),
    @(8, 3470.733333333333, 3720, 3117, 0.05332081317901612),
    @(9, 3068.066666666667, 3354, 2789, 0.05133456389109294)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
